$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "HAWKGUY"
$ws.Range("E2").Select()
